# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the ALC, ARM, BSM, CRP, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value2 = 11905178
$ws.Cells.Item(58, 10).Value2 = 1177.6666
$ws.Cells.Item(58, 12).Value2 = 3532.9998
$ws.Cells.Item(58, 14).Value2 = -3832.9998

$ws.Cells.Item(61, 8).Value2 = 998
$ws.Cells.Item(61, 9).Value2 = 998
$ws.Cells.Item(61, 11).Value2 = 2994
$ws.Cells.Item(61, 13).Value2 = -2822

$ws.Cells.Item(64, 8).Value2 = 8749.950000000001
$ws.Cells.Item(64, 9).Value2 = 3000
$ws.Cells.Item(64, 10).Value2 = 9388.833000000001
$ws.Cells.Item(64, 11).Value2 = 3000
$ws.Cells.Item(64, 12).Value2 = 9388.833000000001
$ws.Cells.Item(64, 13).Value2 = -2752
$ws.Cells.Item(64, 14).Value2 = -9884.833000000001

$ws.Cells.Item(67, 8).Value2 = 8749.950000000001
$ws.Cells.Item(67, 9).Value2 = 3000
$ws.Cells.Item(67, 10).Value2 = 9388.833000000001
$ws.Cells.Item(67, 11).Value2 = 3000
$ws.Cells.Item(67, 12).Value2 = 9388.833000000001
$ws.Cells.Item(67, 13).Value2 = -2142
$ws.Cells.Item(67, 14).Value2 = -11104.833

$ws.Cells.Item(76, 8).Value2 = 4434.091
$ws.Cells.Item(76, 9).Value2 = 4477.5
$ws.Cells.Item(76, 10).Value2 = 4000
$ws.Cells.Item(76, 11).Value2 = 4477.5
$ws.Cells.Item(76, 12).Value2 = 4000
$ws.Cells.Item(76, 13).Value2 = -4162.5
$ws.Cells.Item(76, 14).Value2 = -4630

$ws.Cells.Item(79, 8).Value2 = 4434.091
$ws.Cells.Item(79, 9).Value2 = 4477.5
$ws.Cells.Item(79, 10).Value2 = 4000
$ws.Cells.Item(79, 11).Value2 = 4477.5
$ws.Cells.Item(79, 12).Value2 = 4000
$ws.Cells.Item(79, 13).Value2 = -3385.5
$ws.Cells.Item(79, 14).Value2 = -6184

$ws.Cells.Item(86, 8).Value2 = 4077.4546
$ws.Cells.Item(86, 9).Value2 = 3655.7778
$ws.Cells.Item(86, 11).Value2 = 3655.7778
$ws.Cells.Item(86, 13).Value2 = -2532.7778

$ws.Cells.Item(89, 8).Value2 = 4077.4546
$ws.Cells.Item(89, 9).Value2 = 3655.7778
$ws.Cells.Item(89, 11).Value2 = 18278.889
$ws.Cells.Item(89, 13).Value2 = -12662.889

$ws.Cells.Item(96, 8).Value2 = 2094.4
$ws.Cells.Item(96, 10).Value2 = 4866.6665
$ws.Cells.Item(96, 12).Value2 = 14599.9995
$ws.Cells.Item(96, 14).Value2 = -17345.9995

$ws.Cells.Item(111, 8).Value2 = 1697.45
$ws.Cells.Item(111, 9).Value2 = 1546.0667
$ws.Cells.Item(111, 11).Value2 = 4638.2001
$ws.Cells.Item(111, 13).Value2 = -1571.2001

$ws.Cells.Item(129, 8).Value2 = 1664.125
$ws.Cells.Item(129, 10).Value2 = 2999.3333
$ws.Cells.Item(129, 12).Value2 = 8997.999899999999
$ws.Cells.Item(129, 14).Value2 = -18997.9999

$ws.Cells.Item(137, 8).Value2 = 2071.125
$ws.Cells.Item(137, 9).Value2 = 2179.1667
$ws.Cells.Item(137, 11).Value2 = 6537.500100000001
$ws.Cells.Item(137, 13).Value2 = -3987.500100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(56, 8).Value2 = 44899
$ws.Cells.Item(56, 9).Value2 = 43500
$ws.Cells.Item(56, 10).Value2 = 45831.668
$ws.Cells.Item(56, 11).Value2 = 43500
$ws.Cells.Item(56, 12).Value2 = 45831.668
$ws.Cells.Item(56, 13).Value2 = -42758
$ws.Cells.Item(56, 14).Value2 = -47315.668

$ws.Cells.Item(74, 8).Value2 = 2063.7693
$ws.Cells.Item(74, 9).Value2 = 1982.2
$ws.Cells.Item(74, 11).Value2 = 1982.2
$ws.Cells.Item(74, 13).Value2 = -1108.2

$ws.Cells.Item(77, 8).Value2 = 2063.7693
$ws.Cells.Item(77, 9).Value2 = 1982.2
$ws.Cells.Item(77, 11).Value2 = 9911
$ws.Cells.Item(77, 13).Value2 = -5543

$ws.Cells.Item(97, 8).Value2 = 518.4211
$ws.Cells.Item(97, 9).Value2 = 420.7857
$ws.Cells.Item(97, 10).Value2 = 791.8
$ws.Cells.Item(97, 11).Value2 = 420.7857
$ws.Cells.Item(97, 12).Value2 = 791.8
$ws.Cells.Item(97, 13).Value2 = 75.21429999999998
$ws.Cells.Item(97, 14).Value2 = -1783.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value2 = 2532.1177
$ws.Cells.Item(20, 9).Value2 = 2145
$ws.Cells.Item(20, 11).Value2 = 2145
$ws.Cells.Item(20, 13).Value2 = -1898

$ws.Cells.Item(134, 8).Value2 = 3180.7827
$ws.Cells.Item(134, 9).Value2 = 3041.1875
$ws.Cells.Item(134, 10).Value2 = 3499.8572
$ws.Cells.Item(134, 11).Value2 = 9123.5625
$ws.Cells.Item(134, 12).Value2 = 10499.5716
$ws.Cells.Item(134, 13).Value2 = -6588.5625
$ws.Cells.Item(134, 14).Value2 = -15569.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(69, 8).Value2 = 36052.777
$ws.Cells.Item(69, 9).Value2 = 21895
$ws.Cells.Item(69, 11).Value2 = 21895
$ws.Cells.Item(69, 13).Value2 = -21146

$ws.Cells.Item(72, 8).Value2 = 36052.777
$ws.Cells.Item(72, 9).Value2 = 21895
$ws.Cells.Item(72, 11).Value2 = 65685
$ws.Cells.Item(72, 13).Value2 = -61941

$ws.Cells.Item(81, 8).Value2 = 68664
$ws.Cells.Item(81, 10).Value2 = 68664
$ws.Cells.Item(81, 12).Value2 = 68664
$ws.Cells.Item(81, 14).Value2 = -70660

$ws.Cells.Item(84, 8).Value2 = 68664
$ws.Cells.Item(84, 10).Value2 = 68664
$ws.Cells.Item(84, 12).Value2 = 205992
$ws.Cells.Item(84, 14).Value2 = -215976

$ws.Cells.Item(86, 8).Value2 = 8188.375
$ws.Cells.Item(86, 9).Value2 = 8167.8335
$ws.Cells.Item(86, 11).Value2 = 8167.8335
$ws.Cells.Item(86, 13).Value2 = -7044.8335

$ws.Cells.Item(89, 8).Value2 = 8188.375
$ws.Cells.Item(89, 9).Value2 = 8167.8335
$ws.Cells.Item(89, 11).Value2 = 40839.1675
$ws.Cells.Item(89, 13).Value2 = -35223.1675

$ws.Cells.Item(132, 8).Value2 = 1587.8334
$ws.Cells.Item(132, 9).Value2 = 1614
$ws.Cells.Item(132, 11).Value2 = 4842
$ws.Cells.Item(132, 13).Value2 = -2312

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value2 = 222.89655
$ws.Cells.Item(2, 10).Value2 = 270.2143
$ws.Cells.Item(2, 12).Value2 = 270.2143
$ws.Cells.Item(2, 14).Value2 = -496.2143

$ws.Cells.Item(70, 8).Value2 = 7787.9443
$ws.Cells.Item(70, 9).Value2 = 8147.5
$ws.Cells.Item(70, 10).Value2 = 6529.5
$ws.Cells.Item(70, 11).Value2 = 8147.5
$ws.Cells.Item(70, 12).Value2 = 6529.5
$ws.Cells.Item(70, 13).Value2 = -7877.5
$ws.Cells.Item(70, 14).Value2 = -7069.5

$ws.Cells.Item(73, 8).Value2 = 7787.9443
$ws.Cells.Item(73, 9).Value2 = 8147.5
$ws.Cells.Item(73, 10).Value2 = 6529.5
$ws.Cells.Item(73, 11).Value2 = 8147.5
$ws.Cells.Item(73, 12).Value2 = 6529.5
$ws.Cells.Item(73, 13).Value2 = -7211.5
$ws.Cells.Item(73, 14).Value2 = -8401.5

$ws.Cells.Item(80, 8).Value2 = 5215.273
$ws.Cells.Item(80, 10).Value2 = 5963.933
$ws.Cells.Item(80, 12).Value2 = 5963.933
$ws.Cells.Item(80, 14).Value2 = -7959.933

$ws.Cells.Item(83, 8).Value2 = 5215.273
$ws.Cells.Item(83, 10).Value2 = 5963.933
$ws.Cells.Item(83, 12).Value2 = 29819.665
$ws.Cells.Item(83, 14).Value2 = -39803.665

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value2 = 5149.3794
$ws.Cells.Item(40, 9).Value2 = 3999
$ws.Cells.Item(40, 11).Value2 = 3999
$ws.Cells.Item(40, 13).Value2 = -3863

$ws.Cells.Item(56, 8).Value2 = 7178.5713
$ws.Cells.Item(56, 10).Value2 = 13500
$ws.Cells.Item(56, 12).Value2 = 13500
$ws.Cells.Item(56, 14).Value2 = -14882

$ws.Cells.Item(100, 8).Value2 = 5792.909
$ws.Cells.Item(100, 9).Value2 = 4565
$ws.Cells.Item(100, 10).Value2 = 7566.5557
$ws.Cells.Item(100, 11).Value2 = 4565
$ws.Cells.Item(100, 12).Value2 = 7566.5557
$ws.Cells.Item(100, 13).Value2 = -4024
$ws.Cells.Item(100, 14).Value2 = -8648.555700000001

$ws.Cells.Item(116, 8).Value2 = 214593.8
$ws.Cells.Item(116, 10).Value2 = 214593.8
$ws.Cells.Item(116, 12).Value2 = 214593.8
$ws.Cells.Item(116, 14).Value2 = -223771.8

$ws.Cells.Item(134, 8).Value2 = 19999
$ws.Cells.Item(134, 10).Value2 = 19999
$ws.Cells.Item(134, 12).Value2 = 19999
$ws.Cells.Item(134, 14).Value2 = -30139

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value2 = 16537.223
$ws.Cells.Item(41, 10).Value2 = 16641.857
$ws.Cells.Item(41, 12).Value2 = 16641.857
$ws.Cells.Item(41, 14).Value2 = -17421.857

$ws.Cells.Item(74, 8).Value2 = 7805
$ws.Cells.Item(74, 10).Value2 = 7805
$ws.Cells.Item(74, 12).Value2 = 7805
$ws.Cells.Item(74, 14).Value2 = -9677

$ws.Cells.Item(77, 8).Value2 = 7805
$ws.Cells.Item(77, 10).Value2 = 7805
$ws.Cells.Item(77, 12).Value2 = 23415
$ws.Cells.Item(77, 14).Value2 = -32775

$ws.Cells.Item(126, 8).Value2 = 2009.3334
$ws.Cells.Item(126, 9).Value2 = 1611.5217
$ws.Cells.Item(126, 11).Value2 = 4834.5651
$ws.Cells.Item(126, 13).Value2 = -2364.5651

$ws.Cells.Item(132, 8).Value2 = 3704.0454
$ws.Cells.Item(132, 9).Value2 = 2918.25
$ws.Cells.Item(132, 11).Value2 = 8754.75
$ws.Cells.Item(132, 13).Value2 = -6224.75

$ws.Cells.Item(136, 8).Value2 = 2453.0942
$ws.Cells.Item(136, 9).Value2 = 1305.0294
$ws.Cells.Item(136, 10).Value2 = 4507.5264
$ws.Cells.Item(136, 11).Value2 = 3915.0882
$ws.Cells.Item(136, 12).Value2 = 13522.5792
$ws.Cells.Item(136, 13).Value2 = -1365.0882
$ws.Cells.Item(136, 14).Value2 = -18622.5792
